{"js": "// Find the target paragraph (the last paragraph in the body, ending with\n// \"...m\u0103n\u00e2nc\u0103 multe resurse acum\") and insert four new paragraphs after it:\n//   1. an empty paragraph\n//   2. \"Av\u00e2nd cluster-ul gata, acum revenim la creearea aplicatiilor de blocare a resurselor.\"\n//   3. \"O s\u0103 rescriu aplicatiile din 4 aplica\u021bii separate intr un singur rest api cu flask care, ...\"\n//   4. \"Am modificat aplica\u021bia \u0219i avem un singur server de python care prime\u0219te request-uri care blocheaz\u0103 anumite resurse. Am f\u0103cut \u0219i documentatie cu swagger pentru o testare mai usoara\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"gata avem script automat de creeare a unui cluster yupy, trb pu\u021bin optimizat ca m\u0103n\u00e2nc\u0103 multe resurse acum\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + marker);\n}\n\nconst texts = [\n  \"\",\n  \"Av\u00e2nd cluster-ul gata, acum revenim la creearea aplicatiilor de blocare a resurselor.\",\n  \"O s\u0103 rescriu aplicatiile din 4 aplica\u021bii separate intr un singur rest api cu flask care, \u00een func\u021bie de call-uri, blocheaz\u0103 anumite resurse, \u0219i astfel prin deploy-erea unui singur pod cu aceasta aplica\u021bie pe fiecare nod din cluster, pot altera resursele disponibile de pe fiecare nod. Eliminam \u0219i overhead-ul necesar sa gestionam 4 pod-uri cu 4 aplica\u021bii separate care fiecare sa blocheze o singura resursa.\",\n  \"Am modificat aplica\u021bia \u0219i avem un singur server de python care prime\u0219te request-uri care blocheaz\u0103 anumite resurse. Am f\u0103cut \u0219i documentatie cu swagger pentru o testare mai usoara\",\n];\n\nlet insertAfter = anchor;\nfor (const t of texts) {\n  const p = insertAfter.insertParagraph(t, Word.InsertLocation.after);\n  insertAfter = p;\n}\n\nawait context.sync();\n", "ps1": "# Append four paragraphs at the very end of the document body, right after\n# the paragraph ending in \"...m\u0103n\u00e2nc\u0103 multe resurse acum\":\n#   1. an empty paragraph (blank line)\n#   2. \"Av\u00e2nd cluster-ul gata, acum revenim la creearea aplicatiilor de blocare a resurselor.\"\n#   3. \"O s\u0103 rescriu aplicatiile din 4 aplica\u021bii separate intr un singur rest api cu flask care, ...\"\n#   4. \"Am modificat aplica\u021bia \u0219i avem un singur server de python care prime\u0219te request-uri care blocheaz\u0103 anumite resurse. Am f\u0103cut \u0219i documentatie cu swagger pentru o testare mai usoara\"\n\n$d = $word.ActiveDocument\n\n# --- 1. Blank paragraph, appended at the end of the story ---------------\n$endRange = $d.Range($d.Content.End, $d.Content.End)\n$blankXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"Normal\"/><w:rPr><w:lang w:val=\"ro-RO\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>'\n$endRange.InsertXML($blankXml)\n\n# Re-apply the explicit (redundant-with-style) spacing so it is written out\n# on the paragraph, matching the rest of the document's paragraphs.\n$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$blankPara.Range.ParagraphFormat.SpaceBefore = 0\n$blankPara.Range.ParagraphFormat.SpaceAfter = 8\n\n# --- 2-4. Regular text paragraphs, each appended after the previous one -\n$texts = @(\n  \"Av\u00e2nd cluster-ul gata, acum revenim la creearea aplicatiilor de blocare a resurselor.\",\n  \"O s\u0103 rescriu aplicatiile din 4 aplica\u021bii separate intr un singur rest api cu flask care, \u00een func\u021bie de call-uri, blocheaz\u0103 anumite resurse, \u0219i astfel prin deploy-erea unui singur pod cu aceasta aplica\u021bie pe fiecare nod din cluster, pot altera resursele disponibile de pe fiecare nod. Eliminam \u0219i overhead-ul necesar sa gestionam 4 pod-uri cu 4 aplica\u021bii separate care fiecare sa blocheze o singura resursa.\",\n  \"Am modificat aplica\u021bia \u0219i avem un singur server de python care prime\u0219te request-uri care blocheaz\u0103 anumite resurse. Am f\u0103cut \u0219i documentatie cu swagger pentru o testare mai usoara\"\n)\n\nforeach ($t in $texts) {\n  $tailPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n  $tailRange = $tailPara.Range\n  $tailRange.Collapse(0)\n  $tailRange.InsertParagraphAfter()\n  $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n  $newPara.Range.InsertAfter($t)\n}\n"}
